$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.117.27"
$ws.Range("D2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.00"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  -0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4567"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +7.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3739"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07334"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8612"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.01"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.831.88"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.693"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.00"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.90%  "

$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07072"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.14%  "

$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008841"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.184.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.193"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.005"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.04"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.235"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.66"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.267"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.45"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08882"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.195"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7593"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.981"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.476"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.104"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01971"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05291"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5393"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.215"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.886"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("E42").Value = "  +2.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5231"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +11.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.636"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.71"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.969"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +9.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.30"
$ws.Range("D47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.678"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9999"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.66%  "

$ws.Range("E50").Value = "  -0.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9246"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.65%  "
